# Update the username values in Sheet1 (rows 2-4) from the "_60x" tag
# series to the "_70x" tag series, per "adding tag changes in runner file".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A2 keeps its original rich-text formatting: the trailing digit stays
# rendered in green Courier New (10pt) while the rest of the text changes.
$cellA2 = $ws.Range("A2")
$cellA2.Value = "ds_algo_abs_ch_700"
$lastChar = $cellA2.Characters(18, 1)
$lastChar.Font.Color = 32768
$lastChar.Font.Name = "Courier New"
$lastChar.Font.Size = 10

# A3 and A4 are plain strings - straightforward text replacement.
$ws.Range("A3").Value = "ds_algo_abs_ch_701"
$ws.Range("A4").Value = "ds_algo_abs_ch_702"
